# Ireland.xlsx update: add HPSC data for 2020-04-11 .. 2020-04-16 to the
# "Ireland-manual" sheet (rows 162-173), two rows per date (confirmed/death).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ireland-manual")
$ws.Activate()

$lat = 53.1424
$long = -7.6921
$country = "Ireland"

$rows = @(
    @{ Row = 162; Date = "2020-04-11"; Cases = 553; Type = "confirmed" },
    @{ Row = 163; Date = "2020-04-11"; Cases = 33;  Type = "death" },
    @{ Row = 164; Date = "2020-04-12"; Cases = 430; Type = "confirmed" },
    @{ Row = 165; Date = "2020-04-12"; Cases = 14;  Type = "death" },
    @{ Row = 166; Date = "2020-04-13"; Cases = 527; Type = "confirmed" },
    @{ Row = 167; Date = "2020-04-13"; Cases = 31;  Type = "death" },
    @{ Row = 168; Date = "2020-04-14"; Cases = 548; Type = "confirmed" },
    @{ Row = 169; Date = "2020-04-14"; Cases = 41;  Type = "death" },
    @{ Row = 170; Date = "2020-04-15"; Cases = 657; Type = "confirmed" },
    @{ Row = 171; Date = "2020-04-15"; Cases = 38;  Type = "death" },
    @{ Row = 172; Date = "2020-04-16"; Cases = 629; Type = "confirmed" },
    @{ Row = 173; Date = "2020-04-16"; Cases = 43;  Type = "death" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $country
    $ws.Cells.Item($row, 3).Value = $lat
    $ws.Cells.Item($row, 4).Value = $long
    # Force the date column to text format *before* assigning the value so
    # the "2020-04-xx" string isn't auto-parsed into a date serial (matches
    # rows 162/163, which already carried this text style from the template).
    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $r.Date
    $ws.Cells.Item($row, 6).Value = $r.Cases
    $ws.Cells.Item($row, 7).Value = $r.Type
}

$ws.Range("F172").Select()

# Best-effort: line up the visible viewport with the new rows (saved as
# sheetView/topLeftCell in the OOXML). Not all hosts persist this.
$win = $excel.ActiveWindow
$win.ScrollRow = 137
$win.ScrollColumn = 1
